# Generate Report for Handoff
# The localization-status report refreshes the "Latest Handoff Datetime"
# column (E) for the 678a0f76-d84f-4dc4-9b45-be08ef3ae90e record, which was
# (re-)handed off, on both the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-13 16:43:10"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-13 16:43:13"
